$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range entirely (A1:C6) since rows 5 and 6 and
# several other cells are being removed.
$ws.Range("A1:D6").ClearContents()

# Header row
$ws.Range("A1").Value = "id_wh"
$ws.Range("B1").Value = "serial_number"
$ws.Range("C1").Value = "comment"

# Sample data row
$ws.Range("A2").Value = 17
$ws.Range("B2").Value = "HGST-00001"
$ws.Range("C2").Value = "SAMPLE"
$ws.Range("D2").Value = "id_wh needs to be filled with id of the selected warehouse"

# Explanatory notes column D
$ws.Range("D3").Value = "serial number will auto generate if empty"
$ws.Range("D4").Value = "comment will auto generate if empty"

$ws.Range("D4").Select() | Out-Null

$wb.Save()
